$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: replace "hemintest" / "hemnintest" values with "SagarTest"
$ws.Range("A2").Value = "SagarTest"
$ws.Range("B2").Value = "SagarTest"
$ws.Range("C2").Value = "SagarTest"

# Row 3: keep "testlast" text (re-assert so shared-string table collapses to 5 uniques)
$ws.Range("A3").Value = "testlast"
$ws.Range("B3").Value = "testlast"
$ws.Range("C3").Value = "testlast"

# Active selection ends on C2, matching the author's last edit location
$ws.Range("C2").Select() | Out-Null
